$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: label swaps from "LSPMW" to "LSPM"; B/C get what used to be row 9's values
$ws.Range("A8").Value = "LSPM"
$ws.Range("B8").Value = 109414336465166.3
$ws.Range("C8").Value = 253296144801917.1

# Row 9: label swaps from "LSPM" to "LSPMW"; B/C get new values
$ws.Range("A9").Value = "LSPMW"
$ws.Range("B9").Value = 106412109605695.9
$ws.Range("C9").Value = 246345938679660.8
